$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was updated from 45212 (2023-10-13)
# to 45221 (2023-10-22) for every data row (rows 2 through 158).
$ws.Range("C2:C158").Value = 45221
